$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells, copying the header style from an
# existing header cell (AC1) so they match the bold/centered/bordered look.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Populate the team record (Wins/Losses/Ties) for every player row.
for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 30).Value = 95  # AD
    $ws.Cells.Item($r, 31).Value = 67  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}

Write-Output "done"
